$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily allocation row for 2025-10-19
# Force the date column to be stored as plain text (matching the existing
# "Date" column cells, which are inline/shared strings like "10/18/2025"),
# not auto-converted into an Excel date serial number.
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "10/19/2025"
$ws.Range("A48").ClearFormats()

$ws.Range("B48").Value = 0.1940814947980256
$ws.Range("C48").Value = 0.8059185052019744
